# Update vm_pu results for Case_0_162 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ 2 = 1.02; 3 = 1.054612221764878; 4 = 1.054017034827623; 5 = 1.060312009897332; 6 = 1.069158364782179; 9 = 1.038027765072136; 10 = 1.059623406340892; 11 = 1.056761405992354; 12 = 1.063039132772355; 13 = 1.07186161991257; 14 = 1.0611281931184 }
    3 = @{ 2 = 1.02; 3 = 1.056148000956749; 4 = 1.055386341723682; 5 = 1.061720429266409; 6 = 1.070749730419641; 9 = 1.038306612780032; 10 = 1.060807698869054; 11 = 1.057942450319167; 12 = 1.064260453220431; 13 = 1.073267179786674; 14 = 1.062314167477791 }
    4 = @{ 2 = 1.02; 3 = 1.057140051087728; 4 = 1.056271012671008; 5 = 1.062630511815077; 6 = 1.071778416922958; 9 = 1.038484902330824; 10 = 1.061571921698807; 11 = 1.058704772259468; 12 = 1.065048945061001; 13 = 1.074175146172113; 14 = 1.063079475591621 }
    5 = @{ 2 = 1.02; 3 = 1.057556710015006; 4 = 1.056642608947908; 5 = 1.063012816068391; 6 = 1.072210637305405; 9 = 1.038559343878197; 10 = 1.06189270609367; 11 = 1.059024805463413; 12 = 1.065380006289283; 13 = 1.074556497449573; 14 = 1.063400715537119 }
    6 = @{ 2 = 1.02; 3 = 1.05762664573216; 4 = 1.056704983063163; 5 = 1.063076989569102; 6 = 1.072283195200916; 9 = 1.038571813005814; 10 = 1.061946538368225; 11 = 1.059078514416591; 12 = 1.065435568438637; 13 = 1.074620507236123; 14 = 1.063454624259678 }
    7 = @{ 2 = 1.02; 3 = 1.057145620061875; 4 = 1.056275979205338; 5 = 1.062635621331631; 6 = 1.071784193203757; 9 = 1.038485899029595; 10 = 1.061576209972337; 11 = 1.058709050307355; 12 = 1.065053370362685; 13 = 1.074180243203462; 14 = 1.063083769954991 }
    8 = @{ 2 = 1.02; 3 = 1.055131602978356; 4 = 1.054480084608921; 5 = 1.060788256688461; 6 = 1.069696393062958; 9 = 1.038122447883495; 10 = 1.060024080871037; 11 = 1.057160941082723; 12 = 1.063452257015356; 13 = 1.072336955563506; 14 = 1.061529436652336 }
    9 = @{ 2 = 1.02; 3 = 1.051569234914294; 4 = 1.051304753853303; 5 = 1.057523009411143; 6 = 1.066009147853151; 9 = 1.037465502371153; 10 = 1.057272711201255; 11 = 1.05441820074219; 12 = 1.060616938686928; 13 = 1.069076842241671; 14 = 1.058774159722036 }
    10 = @{ 2 = 1.02; 3 = 1.049184771810456; 4 = 1.049180223701019; 5 = 1.055339045301237; 6 = 1.063544916046651; 9 = 1.037016331502775; 10 = 1.055427090924428; 11 = 1.05257938252992; 12 = 1.058716941486082; 13 = 1.066894902900746; 14 = 1.05692591845272 }
    11 = @{ 2 = 1.02; 3 = 1.048149888265909; 4 = 1.048258370470721; 5 = 1.05439157443122; 6 = 1.062476324086376; 9 = 1.036819151368152; 10 = 1.05462512922351; 11 = 1.051780617260046; 12 = 1.057891810163773; 13 = 1.065947972929644; 14 = 1.056122817874201 }
    12 = @{ 2 = 1.02; 3 = 1.047765116211284; 4 = 1.04791565672117; 5 = 1.05403936255923; 6 = 1.062079156760704; 9 = 1.03674550401138; 10 = 1.054326817031884; 11 = 1.051483530282802; 12 = 1.057584948173515; 13 = 1.065595911145168; 14 = 1.055824082045045 }
    13 = @{ 2 = 1.02; 3 = 1.047847668044047; 4 = 1.047989183531357; 5 = 1.054114925965402; 6 = 1.062164361731681; 9 = 1.036761320023402; 10 = 1.05439082550016; 11 = 1.051547274219455; 12 = 1.057650788049846; 13 = 1.06567144466303; 14 = 1.055888181412689 }
    14 = @{ 2 = 1.02; 3 = 1.048118090500639; 4 = 1.048230047740354; 5 = 1.054362466229487; 6 = 1.062443499160094; 9 = 1.036813071947826; 10 = 1.054600479419551; 11 = 1.051756067958055; 12 = 1.057866452479923; 13 = 1.065918878173187; 14 = 1.056098133064692 }
    15 = @{ 2 = 1.02; 3 = 1.048284657330905; 4 = 1.048378412601479; 5 = 1.054514946727421; 6 = 1.062615452223109; 9 = 1.036844904164675; 10 = 1.054729597180692; 11 = 1.051884660763694; 12 = 1.05799928099112; 13 = 1.06607128615039; 14 = 1.056227434187864 }
    16 = @{ 2 = 1.02; 3 = 1.049253402274482; 4 = 1.049241362953572; 5 = 1.055401887168293; 6 = 1.063615801201883; 9 = 1.037029360907754; 10 = 1.055480254905574; 11 = 1.052632339681194; 12 = 1.058771651075054; 13 = 1.066957701762809; 14 = 1.056979157932817 }
    17 = @{ 2 = 1.02; 3 = 1.04986042180809; 4 = 1.04978214965679; 5 = 1.055957753532236; 6 = 1.064242868232806; 9 = 1.037144344960398; 10 = 1.0559503685906; 11 = 1.053100652199868; 12 = 1.0592554849643; 13 = 1.067513148741204; 14 = 1.057449939233198 }
    18 = @{ 2 = 1.02; 3 = 1.050214255626666; 4 = 1.05009739693104; 5 = 1.056281807887644; 6 = 1.064608475652545; 9 = 1.037211154192392; 10 = 1.056224308715836; 11 = 1.053373565948827; 12 = 1.059537464129246; 13 = 1.067836926097569; 14 = 1.057724268384834 }
    19 = @{ 2 = 1.02; 3 = 1.050334865046154; 4 = 1.050204857067163; 5 = 1.056392272947489; 6 = 1.064733113158657; 9 = 1.037233890548916; 10 = 1.05631766978882; 11 = 1.05346658115523; 12 = 1.059633572461743; 13 = 1.067947291227117; 14 = 1.057817762041251 }
    20 = @{ 2 = 1.02; 3 = 1.049795318263863; 4 = 1.049724147501274; 5 = 1.055898132297697; 6 = 1.064175605478834; 9 = 1.037132035068173; 10 = 1.055899957746252; 11 = 1.053050432097699; 12 = 1.059203598319989; 13 = 1.067453575854505; 14 = 1.05739945679967 }
    21 = @{ 2 = 1.02; 3 = 1.048038468131038; 4 = 1.048159127483757; 5 = 1.054289579562624; 6 = 1.062361306978733; 9 = 1.036797843518751; 10 = 1.054538753441532; 11 = 1.051694594233564; 12 = 1.05780295499897; 13 = 1.065846024348738; 14 = 1.056036319428704 }
    22 = @{ 2 = 1.02; 3 = 1.046931719502506; 4 = 1.047173418153666; 5 = 1.053276600078509; 6 = 1.061219164862826; 9 = 1.036585374598717; 10 = 1.05368043058406; 11 = 1.05083986494685; 12 = 1.056920161730361; 13 = 1.06483338012619; 14 = 1.055176777654325 }
    23 = @{ 2 = 1.02; 3 = 1.04751863490425; 4 = 1.047696127065389; 5 = 1.053813755969205; 6 = 1.061824773967984; 9 = 1.036698231867847; 10 = 1.054135681454307; 11 = 1.051293189928279; 12 = 1.057388353900532; 13 = 1.065370386304187; 14 = 1.055632675033024 }
    24 = @{ 2 = 1.02; 3 = 1.049824736465255; 4 = 1.049750356748863; 5 = 1.055925073102269; 6 = 1.06420599908629; 9 = 1.037137598179032; 10 = 1.055922737068992; 11 = 1.053073125157445; 12 = 1.059227044381119; 13 = 1.067480494912798; 14 = 1.057422268471661 }
    25 = @{ 2 = 1.02; 3 = 1.05249183655064; 4 = 1.052126968696828; 5 = 1.058368378825539; 6 = 1.066963421451337; 9 = 1.037637305390964; 10 = 1.057985981305319; 11 = 1.055129055761859; 12 = 1.061351631257155; 13 = 1.069921129361346; 14 = 1.059488442751462 }
}

foreach ($r in $newValues.Keys) {
    foreach ($c in $newValues[$r].Keys) {
        $ws.Cells.Item($r, $c).Value = $newValues[$r][$c]
    }
}
